$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 1230
$ws.Range("I2").Value = 3158
$ws.Range("J2").Value = 13120
$ws.Range("K2").Value = 57
$ws.Range("L2").Value = 3688
$ws.Range("M2").Value = 215
$ws.Range("N2").Value = 2279
$ws.Range("P2").Value = 55
$ws.Range("Q2").Value = 24
$ws.Range("R2").Value = 186
$ws.Range("S2").Value = 1415
$ws.Range("T2").Value = 2287
$ws.Range("U2").Value = 165
$ws.Range("V2").Value = 20416
$ws.Range("W2").Value = 7
$ws.Range("X2").Value = 20549
$ws.Range("Y2").Value = 25
$ws.Range("Z2").Value = 307
$ws.Range("AA2").Value = 146
